$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-11, columns B-G
$data = @(
    @(0.1374531038001369, 0.9463274094814685, 4.340354840748894, 2.083351828364306, 2.099094096015294, 52),
    @(0.1201387303292092, 0.9941924500794057, 4.368983621713564, 2.090211382064877, 2.10752018100236, 51),
    @(0.1570212086422995, 0.9129929970791744, 4.115459563428113, 2.028659548428004, 2.043107852613974, 50),
    @(0.1668826930103534, 1.018518659023897, 4.584217853051491, 2.141078665778418, 2.156685501776154, 49),
    @(0.1498087523549735, 0.9654930794235904, 4.38482909777947, 2.093998351904669, 2.110735219814999, 48),
    @(0.1482159292065307, 1.037301259063642, 5.28735183827951, 2.299424240604484, 2.327192174260785, 36),
    @(0.1958133101793198, 1.063338432913966, 5.471053399423887, 2.339028302399073, 2.364845904346298, 35),
    @(0.1080685240234339, 1.516053278149146, 9.891401411246379, 3.145059842236134, 3.234328795688273, 18),
    @(-0.8008619424696254, 1.150395842953777, 6.142350263041148, 2.478376537784594, 2.459891756431453, 11),
    @(0.1495946378431339, 0.4169860915672885, 0.2389570720534238, 0.4888323557759078, 0.5203106240290569, 5)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = 2 + $j  # B=2
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
